$d = $word.ActiveDocument

# The document body is a single paragraph that uses <w:br/> line breaks.
# We need to append a batch of new runs right after the very last <w:br/>,
# i.e. at the very end of the document, before the closing </w:p>.
#
# New runs to append (text, italic?):
#   "thuộc "      regular
#   "vào "        regular
#   "kinh "       regular
#   "nghiệm; "    regular
#   "trái "       regular
#   "với "        regular
#   "a "          regular
#   "posteriori. " italic
#   "Suy "        regular
#   "luận "       italic
#   "a "          italic
#   "priort. "    italic

$pos = $d.Content.End

$r = $d.Range($pos, $pos)
$r.InsertAfter("thuộc ")
$pos = $r.End

$r = $d.Range($pos, $pos)
$r.InsertAfter("vào ")
$pos = $r.End

$r = $d.Range($pos, $pos)
$r.InsertAfter("kinh ")
$pos = $r.End

$r = $d.Range($pos, $pos)
$r.InsertAfter("nghiệm; ")
$pos = $r.End

$r = $d.Range($pos, $pos)
$r.InsertAfter("trái ")
$pos = $r.End

$r = $d.Range($pos, $pos)
$r.InsertAfter("với ")
$pos = $r.End

$r = $d.Range($pos, $pos)
$r.InsertAfter("a ")
$pos = $r.End

$insStart = $pos
$r = $d.Range($pos, $pos)
$r.InsertAfter("posteriori. ")
$pos = $r.End
$fmt = $d.Range($insStart, $pos)
$fmt.Font.Italic = 1

$r = $d.Range($pos, $pos)
$r.InsertAfter("Suy ")
$pos = $r.End

$insStart = $pos
$r = $d.Range($pos, $pos)
$r.InsertAfter("luận ")
$pos = $r.End
$fmt = $d.Range($insStart, $pos)
$fmt.Font.Italic = 1

$insStart = $pos
$r = $d.Range($pos, $pos)
$r.InsertAfter("a ")
$pos = $r.End
$fmt = $d.Range($insStart, $pos)
$fmt.Font.Italic = 1

$insStart = $pos
$r = $d.Range($pos, $pos)
$r.InsertAfter("priort. ")
$pos = $r.End
$fmt = $d.Range($insStart, $pos)
$fmt.Font.Italic = 1

Write-Output "done, final pos=$pos"
